$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Simple numeric value updates (row 15-30 data refresh) ---
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 33.333333333333
$ws.Range("I15").Value = 36
$ws.Range("J15").Value = 36
$ws.Range("L15").Value = 44
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 16.129032258064
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 175
$ws.Range("F16").Value = 37
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = 94.736842105263
$ws.Range("I16").Value = 377
$ws.Range("J16").Value = 223
$ws.Range("K16").Value = 69.058295964125
$ws.Range("L16").Value = 51.405622489959
$ws.Range("M16").Value = 12.874251497006
$ws.Range("N16").Value = -73.167259786476
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -46.153846153846
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = -2.380952380952
$ws.Range("I17").Value = 484
$ws.Range("J17").Value = 372
$ws.Range("K17").Value = 30.107526881720
$ws.Range("L17").Value = 48.012232415902
$ws.Range("M17").Value = 109.52380952381
$ws.Range("N17").Value = 7.317073170731
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -65.384615384615
$ws.Range("I18").Value = 145
$ws.Range("J18").Value = 163
$ws.Range("K18").Value = -11.042944785276
$ws.Range("L18").Value = -19.444444444444
$ws.Range("M18").Value = -44.866920152091
$ws.Range("N18").Value = -92.728184553661
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 53
$ws.Range("E19").Value = -73.584905660377
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 126
$ws.Range("H19").Value = -50.793650793650
$ws.Range("I19").Value = 1069
$ws.Range("J19").Value = 634
$ws.Range("K19").Value = 68.611987381703
$ws.Range("L19").Value = 145.183486238532
$ws.Range("M19").Value = 113.8
$ws.Range("N19").Value = -5.565371024734
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 12.5
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 16
$ws.Range("I20").Value = 211
$ws.Range("J20").Value = 178
$ws.Range("K20").Value = 18.539325842696
$ws.Range("L20").Value = 61.068702290076
$ws.Range("M20").Value = 57.462686567164
$ws.Range("N20").Value = -89.836223506743
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 82
$ws.Range("E21").Value = -47.560975609756
$ws.Range("F21").Value = 184
$ws.Range("G21").Value = 241
$ws.Range("H21").Value = -23.651452282157
$ws.Range("I21").Value = 2330
$ws.Range("J21").Value = 1609
$ws.Range("K21").Value = 44.810441267868
$ws.Range("L21").Value = 72.464840858623
$ws.Range("M21").Value = 56.480859637340
$ws.Range("N21").Value = -67.243076057922
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -20
$ws.Range("L22").Value = 192.857142857143
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -32.432432432432
$ws.Range("F24").Value = 186
$ws.Range("H24").Value = 19.230769230769
$ws.Range("I24").Value = 2142
$ws.Range("J24").Value = 1640
$ws.Range("K24").Value = 30.609756097561
$ws.Range("L24").Value = 78.351373855120
$ws.Range("M24").Value = 29.504232164449
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -6.25
$ws.Range("F25").Value = 75
$ws.Range("G25").Value = 62
$ws.Range("H25").Value = 20.967741935483
$ws.Range("I25").Value = 828
$ws.Range("J25").Value = 710
$ws.Range("K25").Value = 16.619718309859
$ws.Range("L25").Value = 30.188679245283
$ws.Range("M25").Value = 43.252595155709
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 54
$ws.Range("J26").Value = 56
$ws.Range("K26").Value = -3.571428571428
$ws.Range("L26").Value = 28.571428571428
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 66.666666666666
$ws.Range("F27").Value = 14
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 180
$ws.Range("I27").Value = 121
$ws.Range("J27").Value = 89
$ws.Range("K27").Value = 35.955056179775
$ws.Range("L27").Value = 59.210526315789
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 75
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = 28.571428571428
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 7
$ws.Range("K30").Value = -36.363636363636
$ws.Range("L30").Value = 40

# --- Cells changing from numeric to the "no data" placeholder text (style 14) ---
# D22 and E22 become shared placeholder strings "0" / "***.*" using style from C22 (s=14)
$ws.Range("D22").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Cells changing from the "no data" placeholder text to real numeric values (rows 28 & 29) ---
# D/G columns use numeric style like D27/G27 (s=15); E/H columns use style like E27/H27 (s=16)
$ws.Range("D27").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("E27").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("G27").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("G28").Value = 1
$ws.Range("H27").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H28").Value = 0
$ws.Range("D27").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("E27").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("G27").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("H27").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = 0
$excel.CutCopyMode = $false
